# Apply "2 Conductor Area Calculations" edit to the "in" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update header row: column C header text "Lb Per Foot" -> "Weight" ---
$ws.Cells.Item(1, 3).Value = "Weight"

# --- 2. Update column C data values (rows 2-12): lb/foot -> weight (x1000, rounded) ---
$ws.Cells.Item(2, 3).Value = 960
$ws.Cells.Item(3, 3).Value = 735
$ws.Cells.Item(4, 3).Value = 743
$ws.Cells.Item(5, 3).Value = 743
$ws.Cells.Item(6, 3).Value = 495
$ws.Cells.Item(7, 3).Value = 711
$ws.Cells.Item(8, 3).Value = 767
$ws.Cells.Item(9, 3).Value = 865
$ws.Cells.Item(10, 3).Value = 1284
$ws.Cells.Item(11, 3).Value = 880
$ws.Cells.Item(12, 3).Value = 860

# --- 3. Replace old rows 13 & 14 ("CABLE 1" / "CABLE 2") ---
# Clear the old contents of B13/C13 and B14/C14 first.
$ws.Cells.Item(13, 2).ClearContents()
$ws.Cells.Item(13, 3).ClearContents()
$ws.Cells.Item(14, 2).ClearContents()
$ws.Cells.Item(14, 3).ClearContents()

# Row 13: single banner cell
$ws.Cells.Item(13, 1).Value = "* 2 Conductor Cables Below *"

# Row 14: new sub-header row (Size / Length / Width / Weight), bold like row 1
$ws.Cells.Item(14, 1).Value = "Size"
$ws.Cells.Item(14, 2).Value = "Length"
$ws.Cells.Item(14, 3).Value = "Width"
$ws.Cells.Item(14, 4).Value = "Weight"
$ws.Range($ws.Cells.Item(14, 1), $ws.Cells.Item(14, 4)).Font.Bold = $true

# --- 4. New data rows 15-18 for the 2-conductor cables ---
$ws.Cells.Item(15, 1).Value = "2C#4"
$ws.Cells.Item(15, 2).Value = 0.8
$ws.Cells.Item(15, 3).Value = 1.284
$ws.Cells.Item(15, 4).Value = 635

$ws.Cells.Item(16, 1).Value = "2C#6"
$ws.Cells.Item(16, 2).Value = 0.74
$ws.Cells.Item(16, 3).Value = 1.1140000000000001
$ws.Cells.Item(16, 4).Value = 540

$ws.Cells.Item(17, 1).Value = "2C#9"
$ws.Cells.Item(17, 2).Value = 0.55100000000000005
$ws.Cells.Item(17, 3).Value = 0.83799999999999997
$ws.Cells.Item(17, 4).Value = 295

$ws.Cells.Item(18, 1).Value = "2C#14"
$ws.Cells.Item(18, 2).Value = 0.48699999999999999
$ws.Cells.Item(18, 3).Value = 0.71699999999999997
$ws.Cells.Item(18, 4).Value = 201

# Center-align the "Length" values for the new cable rows (B15:B18)
$ws.Range($ws.Cells.Item(15, 2), $ws.Cells.Item(18, 2)).HorizontalAlignment = -4108

# --- 5. Column widths (match the widened/added columns from the author's edit) ---
# Column A widened to fit the new "* 2 Conductor Cables Below *" banner text.
$ws.Columns.Item(1).ColumnWidth = 26.0005
# Column C (Weight) / D (Width) / E (Weight of 2nd table) sized to fit their content.
$ws.Columns.Item(3).ColumnWidth = 10.33305
$ws.Columns.Item(4).ColumnWidth = 10.0006
$ws.Columns.Item(5).ColumnWidth = 14.33075

# --- 6. Sheet view: zoom + selected cell ---
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("A13").Select() | Out-Null
